$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Insert 2 new blank rows at row 5, pushing the existing
#     rows 5-7 (Botafogo SP/Avai, Progreso/Fenix, CA Cerro/Boston River)
#     down to rows 7-9. ---
$ws.Rows.Item(5).Resize(2).Insert()

# --- Step 2: Small data correction on row 2 (Banfield vs Tigre). ---
$ws.Range("X2").Value = 10
$ws.Range("Z2").Value = 23
$ws.Range("AK2").Value = 34

# --- Step 3: Odds correction on the match that is now row 7
#     (Botafogo SP vs Avai, BRAZIL - SERIE B). ---
$ws.Range("G7").Value = 2.55
$ws.Range("K7").Value = 1.8
$ws.Range("W7").Value = 6
$ws.Range("AW7").Value = 4.75
$ws.Range("AZ7").Value = 67

# --- Step 4: Fill in the two brand-new matches inserted as rows 5 and 6,
#     plus the additional new match appended as row 10. ---
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY","AZ","BA","BB","BC","BD")

$row5 = @("K4udGliF", "18/11/2024", "21:30", "ARGENTINA - TORNEO BETANO", "Atl. Tucuman", "Huracan", 2.3, 3, 3.5, 3.2, 1.83, 4.5, 1.13, 6, 1.57, 2.25, 2.88, 1.4, 1.67, 2.1, 2.25, 1.57, 5.5, 9, 10, 21, 23, 41, 5.5, 6, 21, 81, 7, 15, 13, 41, 41, 51, 201, 4, 15, 34, 51, 101, 351, 2.1, 10, 81, 5, 23, 41, 81, 151, 451, 126, 126)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value = $row5[$i]
}

$row6 = @("W2Rn64T7", "18/11/2024", "21:30", "ARGENTINA - TORNEO BETANO", "Instituto", "Argentinos Jrs", 2.25, 3, 3.6, 3.1, 1.91, 4.33, 1.11, 6.5, 1.5, 2.5, 2.6, 1.48, 1.57, 2.25, 2.1, 1.67, 6, 9.5, 10, 21, 23, 41, 6, 6, 19, 67, 8, 15, 13, 41, 34, 41, 201, 4, 13, 29, 51, 81, 301, 2.25, 9.5, 81, 5, 21, 34, 81, 126, 351, 126, 126)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "6").Value = $row6[$i]
}

$row10 = @("ADYaA6BG", "18/11/2024", "21:30", "URUGUAY - PRIMERA DIVISION", "Cerro Largo", "Wanderers", 2, 3.1, 4.1, 2.75, 2, 4.5, 1.08, 8, 1.36, 3, 2.25, 1.58, 1.5, 2.5, 2, 1.73, 6, 8.5, 9, 17, 19, 34, 7.5, 6, 17, 51, 10, 19, 15, 41, 41, 41, 1000, 4, 11, 26, 41, 67, 201, 2.5, 9, 67, 5.5, 23, 34, 81, 126, 301, 51, 51)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "10").Value = $row10[$i]
}
